# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.097.30"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").Value = "2.468.35"
$ws.Range("E3").Value = "  -2.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "2.468.09"
$ws.Range("E9").Value = "  -2.91%  "

$ws.Range("E10").Value = "  -2.92%  "

$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("E12").Value = "  -2.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.39%  "

$ws.Range("E15").Value = "  -2.84%  "

$ws.Range("D16").Value = "66.756.91"
$ws.Range("E16").Value = "  -1.68%  "

$ws.Range("E17").Value = "  -4.86%  "

$ws.Range("D18").Value = "2.461.79"
$ws.Range("E18").Value = "  -2.89%  "

$ws.Range("E19").Value = "  -5.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "354.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.17%  "

$ws.Range("E22").Value = "  -2.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.75%  "

$ws.Range("E25").Value = "  -7.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -57.79%  "

$ws.Range("D29").Value = "2.593.06"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").Value = "0.0₃0906"
$ws.Range("E30").Value = "  -6.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "514.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.76"
$ws.Range("D32").Style = "Normal"

$ws.Range("E33").Value = "  -6.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  -9.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("E38").Value = "  +0.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.16%  "

$ws.Range("E40").Value = "  -6.41%  "

$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.58%  "

$ws.Range("E44").Value = "  -7.05%  "

$ws.Range("E45").Value = "  -7.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.93%  "

$ws.Range("E49").Value = "  -6.94%  "

$ws.Range("D50").Value = "0.0₆0256"
$ws.Range("E50").Value = "  -11.46%  "

$ws.Range("E51").Value = "  -7.24%  "
